$d = $word.ActiveDocument

# --- Paragraph 1: "GIT CMD" -> "GIT CMD basicos " with bookmark at end ---
$p1 = $d.Paragraphs(1).Range
$p1.InsertAfter(" basicos X")
$bm = $d.Range(16, 16)
$d.Bookmarks.Add("_GoBack", $bm)
$placeholder = $d.Range(16, 17)
$placeholder.Delete()
